# Update consumption data: shift dates forward by 7 days (03-04 Apr 2025 -> 10-11 Apr 2025)
# and refresh "Actual Consumption (MW)" values with the latest fetched data
# (Improved Wind Production Forecast model).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseDate = Get-Date -Year 2025 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0

$bVals = @(5961,6001,5942,5915,5803,5861,5818,5756,5762,5848,5779,5829,5844,5856,5804,5883,5971,6022,6111,6199,6369,6473,6614,6731,7016,7118,7284,7386,7483,7491,7553,7598,7537,7510,7451,7390,7311,7261,7255,7248,7092,7092,7054,7008,6794,6741,6713,6715,6666,6719,6761,6707,6727,6707,6625,6570,6688,6670,6742,6664,6676,6716,6755,6843,6831,6888,6891,6960,7034,7142,7223,7372,7519,7602,7733,7840,8020,8123,8213,8141,8129,7980,7899,7780,7561,7412,7256,7136,6919,6731,6673,6558,6460,6371,6329,6241,6240,6138,6107,6062,6033,5991,6040,6024,5987,6006,5961,5981,5995,6026,6051,6101,6113,6204,6245,6361,6525,6628,6773,6911,7122,7330,7470,7562,7616,7690,7721,7617,7533,7608,7578,7476,7382,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = 2 + $i
    $dt = $baseDate.AddMinutes(15 * $i)
    $quarter = ($i % 96) + 1
    $lookup = $dt.ToString("dd.MM.yyyy") + $quarter

    $ws.Cells.Item($row, 1).Value = $dt
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
    $ws.Cells.Item($row, 4).Value = $lookup
}
